$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previously used range entirely before rewriting it
$ws.Range("A1:I4").Clear()

# Update header row
$ws.Range("A1").Value = "Qtd_Nós"
$ws.Range("B1").Value = "Ativos"
$ws.Range("C1").Value = "Distancia"
$ws.Range("D1").Value = "Tempo"

# Update data row
$ws.Range("A2").Value = 42
$ws.Range("B2").Value = 30
$ws.Range("C2").Value = 7159
$ws.Range("D2").Value = 2671.142466068268
